# Generate Report for Handoff
# Insert a new "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f" row ahead of the existing
# "b98b5258-999e-49a0-b209-432c0300c06e" row on each of the three sheets
# (Overview, zh-cn, de-de), pushing the b98b5258 row down by one.

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A (File Name) / B (zh-cn) / C (de-de) / D (Latest Handoff Date)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Shift row 6 (b98b5258...) down to row 7, copying formats along the way.
$ws1.Rows.Item(6).Insert()

# Populate the freshly inserted row 6 with the new file's summary data.
$ws1.Range("A6").Value = "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md"
$ws1.Range("B6").Value = "Ready for handoff"
$ws1.Range("C6").Value = "Ready for handoff"
$ws1.Range("D6").Value = "2016-24-11 08:24:48"

# Rebuild every hyperlink on this sheet in final, correct order (the insert
# above does not itself relocate hyperlink anchors).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.md", $missing, $missing, "13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb507c6b1a54eed1a6867280650b51f9e0e93a7/e2e/463f88a5-117f-44f0-adbc-81e9ccadac3c.md", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/845ebbd2-cb09-4904-ac07-1994f9d528bb.md", $missing, $missing, "845ebbd2-cb09-4904-ac07-1994f9d528bb.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1101e5786529ec1c574ec91200f59d0b768a8b4f/e2e/d99be125-a340-45f3-a35c-c9f86371d7c9.md", $missing, $missing, "d99be125-a340-45f3-a35c-c9f86371d7c9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a60a7a7832f60cdcee3691e42ac6c842a7127553/e2e/3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md", $missing, $missing, "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/372183d9726044da73322b109443bb51722b0cee/e2e/b98b5258-999e-49a0-b209-432c0300c06e.md", $missing, $missing, "b98b5258-999e-49a0-b209-432c0300c06e.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": per-locale handoff detail table
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(6).Insert()

$ws2.Range("A6").Value = "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md"
$ws2.Range("B6").Value = ".md"
$ws2.Range("C6").Value = "Ready for handoff"
$ws2.Range("D6").Value = "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.a60a7a7832f60cdcee3691e42ac6c842a7127553.zh-cn.xlf"
$ws2.Range("E6").Value = "2016-03-11 08:24:38"
$ws2.Range("H6").Value = "0001-01-01 00:00:00"
$ws2.Range("I6").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.md", $missing, $missing, "13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.md", $missing, $missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10456c96dc3dad3a64f59ef834662fe9b53d8176/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.aa6cdad381deaf236804fe07ada9e095e481f549.zh-cn.xlf", $missing, $missing, "13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.aa6cdad381deaf236804fe07ada9e095e481f549.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb507c6b1a54eed1a6867280650b51f9e0e93a7/e2e/463f88a5-117f-44f0-adbc-81e9ccadac3c.md", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb507c6b1a54eed1a6867280650b51f9e0e93a7/e2e/463f88a5-117f-44f0-adbc-81e9ccadac3c.md", $missing, $missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a2409c336af79324645e939eb3babce2b2e94991/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/463f88a5-117f-44f0-adbc-81e9ccadac3c.463acbc0e398c2251bf311d81094020db05969f7.zh-cn.xlf", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.463acbc0e398c2251bf311d81094020db05969f7.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/de3f06cb7d55ebbb4845a57b164f77e04106c6d5/e2e/463f88a5-117f-44f0-adbc-81e9ccadac3c.md", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/32f8f35c60f69b07a1669f8a8a03788ff1e54bc0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/463f88a5-117f-44f0-adbc-81e9ccadac3c.463acbc0e398c2251bf311d81094020db05969f7.zh-cn.xlf", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.463acbc0e398c2251bf311d81094020db05969f7.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/845ebbd2-cb09-4904-ac07-1994f9d528bb.md", $missing, $missing, "845ebbd2-cb09-4904-ac07-1994f9d528bb.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/845ebbd2-cb09-4904-ac07-1994f9d528bb.md", $missing, $missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10456c96dc3dad3a64f59ef834662fe9b53d8176/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/845ebbd2-cb09-4904-ac07-1994f9d528bb.0bce75f6262564538c5286fff95b38b4ef05bef3.zh-cn.xlf", $missing, $missing, "845ebbd2-cb09-4904-ac07-1994f9d528bb.0bce75f6262564538c5286fff95b38b4ef05bef3.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1101e5786529ec1c574ec91200f59d0b768a8b4f/e2e/d99be125-a340-45f3-a35c-c9f86371d7c9.md", $missing, $missing, "d99be125-a340-45f3-a35c-c9f86371d7c9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/1101e5786529ec1c574ec91200f59d0b768a8b4f/e2e/d99be125-a340-45f3-a35c-c9f86371d7c9.md", $missing, $missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/166116d06a7e70a6b29b0ee8c49d05bdefb53eb9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/d99be125-a340-45f3-a35c-c9f86371d7c9.111b5e45080c1cd1be26c3526f075e2fb4f7473b.zh-cn.xlf", $missing, $missing, "d99be125-a340-45f3-a35c-c9f86371d7c9.111b5e45080c1cd1be26c3526f075e2fb4f7473b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a60a7a7832f60cdcee3691e42ac6c842a7127553/e2e/3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md", $missing, $missing, "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/a60a7a7832f60cdcee3691e42ac6c842a7127553/e2e/3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md", $missing, $missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a60a7a7832f60cdcee3691e42ac6c842a7127553/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.a60a7a7832f60cdcee3691e42ac6c842a7127553.zh-cn.xlf", $missing, $missing, "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.a60a7a7832f60cdcee3691e42ac6c842a7127553.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/372183d9726044da73322b109443bb51722b0cee/e2e/b98b5258-999e-49a0-b209-432c0300c06e.md", $missing, $missing, "b98b5258-999e-49a0-b209-432c0300c06e.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/372183d9726044da73322b109443bb51722b0cee/e2e/b98b5258-999e-49a0-b209-432c0300c06e.md", $missing, $missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/00f4541a25408fd545cbda9ae1a81e6912e97504/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/b98b5258-999e-49a0-b209-432c0300c06e.2b18ae0ac6123e526871b212279646d45a4cbb3d.zh-cn.xlf", $missing, $missing, "b98b5258-999e-49a0-b209-432c0300c06e.2b18ae0ac6123e526871b212279646d45a4cbb3d.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": per-locale handoff detail table
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(6).Insert()

$ws3.Range("A6").Value = "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md"
$ws3.Range("B6").Value = ".md"
$ws3.Range("C6").Value = "Ready for handoff"
$ws3.Range("D6").Value = "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.a60a7a7832f60cdcee3691e42ac6c842a7127553.de-de.xlf"
$ws3.Range("E6").Value = "2016-03-11 08:24:48"
$ws3.Range("H6").Value = "0001-01-01 00:00:00"
$ws3.Range("I6").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.md", $missing, $missing, "13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.md", $missing, $missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a05d0631cc2a3bede820f81317837e27a4b5723/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.aa6cdad381deaf236804fe07ada9e095e481f549.de-de.xlf", $missing, $missing, "13ba6f46-54e4-4469-8f1c-27fdf8a8d7e6.aa6cdad381deaf236804fe07ada9e095e481f549.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb507c6b1a54eed1a6867280650b51f9e0e93a7/e2e/463f88a5-117f-44f0-adbc-81e9ccadac3c.md", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/5cb507c6b1a54eed1a6867280650b51f9e0e93a7/e2e/463f88a5-117f-44f0-adbc-81e9ccadac3c.md", $missing, $missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f16a1260b2f569d31498d0e3e7e60efc8bd7fa16/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/463f88a5-117f-44f0-adbc-81e9ccadac3c.463acbc0e398c2251bf311d81094020db05969f7.de-de.xlf", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.463acbc0e398c2251bf311d81094020db05969f7.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a3880c5bd36d818ba7811abdf2f678ebe54c2328/e2e/463f88a5-117f-44f0-adbc-81e9ccadac3c.md", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e101fd998c04ad7f9a7a1551d8dbba0cfe49dd0a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/463f88a5-117f-44f0-adbc-81e9ccadac3c.463acbc0e398c2251bf311d81094020db05969f7.de-de.xlf", $missing, $missing, "463f88a5-117f-44f0-adbc-81e9ccadac3c.463acbc0e398c2251bf311d81094020db05969f7.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/845ebbd2-cb09-4904-ac07-1994f9d528bb.md", $missing, $missing, "845ebbd2-cb09-4904-ac07-1994f9d528bb.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/5c95f7474dce86c1c41a8900734ecbc8bdd919b7/e2e/845ebbd2-cb09-4904-ac07-1994f9d528bb.md", $missing, $missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a05d0631cc2a3bede820f81317837e27a4b5723/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/845ebbd2-cb09-4904-ac07-1994f9d528bb.0bce75f6262564538c5286fff95b38b4ef05bef3.de-de.xlf", $missing, $missing, "845ebbd2-cb09-4904-ac07-1994f9d528bb.0bce75f6262564538c5286fff95b38b4ef05bef3.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/1101e5786529ec1c574ec91200f59d0b768a8b4f/e2e/d99be125-a340-45f3-a35c-c9f86371d7c9.md", $missing, $missing, "d99be125-a340-45f3-a35c-c9f86371d7c9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/1101e5786529ec1c574ec91200f59d0b768a8b4f/e2e/d99be125-a340-45f3-a35c-c9f86371d7c9.md", $missing, $missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8ca6474cbb9ac013d46e82e9b91288920eb33989/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/d99be125-a340-45f3-a35c-c9f86371d7c9.111b5e45080c1cd1be26c3526f075e2fb4f7473b.de-de.xlf", $missing, $missing, "d99be125-a340-45f3-a35c-c9f86371d7c9.111b5e45080c1cd1be26c3526f075e2fb4f7473b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a60a7a7832f60cdcee3691e42ac6c842a7127553/e2e/3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md", $missing, $missing, "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/a60a7a7832f60cdcee3691e42ac6c842a7127553/e2e/3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.md", $missing, $missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a60a7a7832f60cdcee3691e42ac6c842a7127553/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.a60a7a7832f60cdcee3691e42ac6c842a7127553.de-de.xlf", $missing, $missing, "3dbb7bb8-5ab7-4748-9001-4d5408a78c9f.a60a7a7832f60cdcee3691e42ac6c842a7127553.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/372183d9726044da73322b109443bb51722b0cee/e2e/b98b5258-999e-49a0-b209-432c0300c06e.md", $missing, $missing, "b98b5258-999e-49a0-b209-432c0300c06e.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/372183d9726044da73322b109443bb51722b0cee/e2e/b98b5258-999e-49a0-b209-432c0300c06e.md", $missing, $missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d1dd7c3b13f90805fb0fbff02ff4dcc57e4908e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/b98b5258-999e-49a0-b209-432c0300c06e.2b18ae0ac6123e526871b212279646d45a4cbb3d.de-de.xlf", $missing, $missing, "b98b5258-999e-49a0-b209-432c0300c06e.2b18ae0ac6123e526871b212279646d45a4cbb3d.de-de.xlf") | Out-Null
